# Apply updated odds values to Sheet1 based on the authoritative diff.
# Each entry maps a cell address (A1 notation) to its new numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "G2" = 2.3
    "I2" = 3
    "Z2" = 15
    "AF2" = 19
    "G7" = 2.55
    "I7" = 2.75
    "N7" = 1.84
    "O7" = 2.06
    "R7" = 1.62
    "S7" = 2.2
    "T7" = 9.5
    "W7" = 23
    "Z7" = 12
    "AH7" = 29
    "AJ7" = 29
    "G8" = 1.67
    "I8" = 5.25
    "L8" = 1.3
    "M8" = 3.5
    "N8" = 2.03
    "O8" = 1.87
    "G9" = 2.7
    "H9" = 2.88
    "I9" = 2.88
    "J9" = 1.13
    "K9" = 6
    "P9" = 1.62
    "Q9" = 2.2
    "T9" = 6.5
    "U9" = 11
    "V9" = 11
    "W9" = 29
    "X9" = 26
    "Z9" = 6
    "AF9" = 12
    "AH9" = 29
    "AI9" = 29
    "G10" = 1.91
    "H10" = 3.25
    "I10" = 4.2
    "AE10" = 8.5
    "G12" = 2.2
    "H12" = 2.92
    "I12" = 3.55
    "N12" = 2.4
    "O12" = 1.5
    "P12" = 1.5
    "Q12" = 2.42
    "S12" = 1.7
    "T12" = 5.7
    "V12" = 9.25
    "W12" = 21
    "X12" = 22
    "Y12" = 40
    "AA12" = 5.7
    "AB12" = 16.5
    "AE12" = 8.5
    "AF12" = 18
    "AG12" = 12.5
    "AI12" = 37
    "AJ12" = 50
    "G16" = 2.7
    "H16" = 2.9
    "I16" = 2.75
    "J16" = 1.13
    "K16" = 6
    "N16" = 2.88
    "O16" = 1.4
    "T16" = 6.5
    "U16" = 12
    "V16" = 12
    "W16" = 29
    "P17" = 1.75
    "Q17" = 2.05
    "G18" = 1.75
    "H18" = 3.1
    "I18" = 5.2
    "L18" = 1.5
    "M18" = 2.25
    "N18" = 2.45
    "T18" = 4.7
    "U18" = 6.6
    "W18" = 13.5
    "Z18" = 6.2
    "AA18" = 6.5
    "AE18" = 10
    "AF18" = 28
    "AI18" = 70
    "J19" = 1.08
    "K19" = 8
    "G99" = 2.3
    "I99" = 3.25
    "L99" = 1.5
    "M99" = 2.5
    "U99" = 10
    "X99" = 23
    "AE99" = 7.5
    "AH99" = 34
    "J100" = 1.08
    "K100" = 8
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
